# Add a new package entry row (row 9) to Sheet1, mirroring the existing
# "DTDemo" row (row 8) but for a ValueMapping artifact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = "DTDemo"
$ws.Range("B9").Value = "DTDemo"
$ws.Range("C9").Value = "VM_SourceAgency1_SourceSchme1_TargetAgency1_TargetSchema1"
$ws.Range("D9").Value = "1.0.0"
$ws.Range("E9").Value = "ValueMapping"

# Row 8's "Date Uploaded" cell (F8) already holds the literal text
# "2026-02-04". Copy it down into F9 instead of typing the string again,
# so the new cell keeps the exact same text representation (rather than
# letting Excel auto-convert the string into a numeric date serial).
$ws.Range("F8").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4163)  # xlPasteValues
